$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.470.46'
$ws.Range("E2").Value = '  +3.51%  '
$ws.Range("D3").Value = '1.591.57'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.06'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.492'
$ws.Range("E7").Value = '  +0.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.40'
$ws.Range("E8").Value = '  +7.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.252'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0601'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("E11").Value = '  +1.66%  '
$ws.Range("D12").Value = '1.818.63'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '1.593.25'
$ws.Range("E13").Value = '  +1.58%  '
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '28.482.05'
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.18'
$ws.Range("E17").Value = '  +1.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.31'
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("E20").Value = '  -0.80%  '
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.70'
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.57'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.107'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").Value = '1.404.73'
$ws.Range("E34").Value = '  -3.64%  '
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("E36").Value = '  -10.14%  '
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.62'
$ws.Range("E38").Value = '  +10.37%  '
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.542'
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.63'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.20'
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").Value = '1.727.64'
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("E48").Value = '  +1.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.31'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("E51").Value = '  -0.70%  '
